$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 199.33333
$ws.Range("I38").Value = 127.71429
$ws.Range("J38").Value = 450
$ws.Range("K38").Value = 383.14287
$ws.Range("L38").Value = 1350
$ws.Range("M38").Value = -11.14287000000002
$ws.Range("N38").Value = -2094
$ws.Range("H43").Value = 2100
$ws.Range("J43").Value = 2100
$ws.Range("L43").Value = 2100
$ws.Range("N43").Value = -2238
$ws.Range("H51").Value = 5958170
$ws.Range("J51").Value = 10211081
$ws.Range("L51").Value = 10211081
$ws.Range("N51").Value = -10212049
$ws.Range("H100").Value = 2163.6365
$ws.Range("J100").Value = 3668
$ws.Range("L100").Value = 3668
$ws.Range("N100").Value = -4750
$ws.Range("H107").Value = 436.16666
$ws.Range("I107").Value = 436.16666
$ws.Range("K107").Value = 436.16666
$ws.Range("M107").Value = 1483.83334
$ws.Range("H111").Value = 32927.43
$ws.Range("I111").Value = 1357.25
$ws.Range("J111").Value = 75021
$ws.Range("K111").Value = 4071.75
$ws.Range("L111").Value = 225063
$ws.Range("M111").Value = -1004.75
$ws.Range("N111").Value = -231197
$ws.Range("H112").Value = 1956.0857
$ws.Range("J112").Value = 1807.7354
$ws.Range("L112").Value = 5423.206200000001
$ws.Range("N112").Value = -7639.206200000001
$ws.Range("H113").Value = 142861280
$ws.Range("J113").Value = 6002
$ws.Range("L113").Value = 6002
$ws.Range("N113").Value = -12510
$ws.Range("H115").Value = 513.3333
$ws.Range("I115").Value = 513.3333
$ws.Range("K115").Value = 1539.9999
$ws.Range("M115").Value = 27.00009999999997
$ws.Range("H132").Value = 4956.4443
$ws.Range("I132").Value = 4404.613
$ws.Range("J132").Value = 8377.799999999999
$ws.Range("K132").Value = 13213.839
$ws.Range("L132").Value = 25133.4
$ws.Range("M132").Value = -10683.839
$ws.Range("N132").Value = -30193.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 15152917
$ws.Range("I61").Value = 16667888
$ws.Range("K61").Value = 16667888
$ws.Range("M61").Value = -16667676
$ws.Range("H74").Value = 2491.5625
$ws.Range("I74").Value = 2133.6428
$ws.Range("J74").Value = 4997
$ws.Range("K74").Value = 2133.6428
$ws.Range("L74").Value = 4997
$ws.Range("M74").Value = -1259.6428
$ws.Range("N74").Value = -6745
$ws.Range("H77").Value = 2491.5625
$ws.Range("I77").Value = 2133.6428
$ws.Range("J77").Value = 4997
$ws.Range("K77").Value = 10668.214
$ws.Range("L77").Value = 24985
$ws.Range("M77").Value = -6300.214
$ws.Range("N77").Value = -33721
$ws.Range("H97").Value = 1156.0834
$ws.Range("I97").Value = 1187.5454
$ws.Range("K97").Value = 1187.5454
$ws.Range("M97").Value = -691.5454
$ws.Range("H132").Value = 45456616
$ws.Range("I132").Value = 55557590
$ws.Range("K132").Value = 166672770
$ws.Range("M132").Value = -166670240
$ws.Range("H136").Value = 15152917
$ws.Range("I136").Value = 16667888
$ws.Range("K136").Value = 50003664
$ws.Range("M136").Value = -50001114

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 7443.6665
$ws.Range("I94").Value = 9165.5
$ws.Range("K94").Value = 9165.5
$ws.Range("M94").Value = -8714.5
$ws.Range("H134").Value = 1384.1666
$ws.Range("I134").Value = 1384.1666
$ws.Range("K134").Value = 4152.4998
$ws.Range("M134").Value = -1617.4998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 40.666668
$ws.Range("I7").Value = 32.6
$ws.Range("J7").Value = 81
$ws.Range("K7").Value = 32.6
$ws.Range("L7").Value = 81
$ws.Range("M7").Value = 80.40000000000001
$ws.Range("N7").Value = -307
$ws.Range("H31").Value = 3123.7144
$ws.Range("I31").Value = 2573.5
$ws.Range("K31").Value = 2573.5
$ws.Range("M31").Value = -2278.5
$ws.Range("H34").Value = 3123.7144
$ws.Range("I34").Value = 2573.5
$ws.Range("K34").Value = 2573.5
$ws.Range("M34").Value = -2371.5
$ws.Range("H99").Value = 2442.5833
$ws.Range("I99").Value = 2231.1
$ws.Range("K99").Value = 2231.1
$ws.Range("M99").Value = -733.0999999999999
$ws.Range("H126").Value = 2442.5833
$ws.Range("I126").Value = 2231.1
$ws.Range("K126").Value = 6693.299999999999
$ws.Range("M126").Value = -4223.299999999999
$ws.Range("H132").Value = 1975.5834
$ws.Range("I132").Value = 1651
$ws.Range("K132").Value = 4953
$ws.Range("M132").Value = -2423
$ws.Range("H134").Value = 2851.4546
$ws.Range("I134").Value = 2536.7
$ws.Range("K134").Value = 7610.099999999999
$ws.Range("M134").Value = -5075.099999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1365
$ws.Range("I107").Value = 1577.6666
$ws.Range("J107").Value = 89
$ws.Range("K107").Value = 4732.9998
$ws.Range("L107").Value = 267
$ws.Range("M107").Value = -2812.9998
$ws.Range("N107").Value = -4107
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("M108").ClearContents()
$ws.Range("H115").Value = 6367.857
$ws.Range("I115").Value = 2166
$ws.Range("J115").Value = 6691.077
$ws.Range("K115").Value = 6498
$ws.Range("L115").Value = 20073.231
$ws.Range("M115").Value = -5323
$ws.Range("N115").Value = -22423.231
$ws.Range("H118").Value = 190
$ws.Range("I118").Value = 190
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 570
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = 673
$ws.Range("N118").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2309.6667
$ws.Range("I80").Value = 2226.0908
$ws.Range("K80").Value = 2226.0908
$ws.Range("M80").Value = -1228.0908
$ws.Range("H83").Value = 2309.6667
$ws.Range("I83").Value = 2226.0908
$ws.Range("K83").Value = 11130.454
$ws.Range("M83").Value = -6138.454
$ws.Range("H97").Value = 2378.3684
$ws.Range("I97").Value = 2386.8462
$ws.Range("J97").Value = 2360
$ws.Range("K97").Value = 2386.8462
$ws.Range("L97").Value = 2360
$ws.Range("M97").Value = -1890.8462
$ws.Range("N97").Value = -3352
$ws.Range("H126").Value = 6316.643
$ws.Range("I126").Value = 8608.933999999999
$ws.Range("K126").Value = 25826.802
$ws.Range("M126").Value = -23356.802

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 26319412
$ws.Range("I7").Value = 45457416
$ws.Range("J7").Value = 4657.375
$ws.Range("K7").Value = 45457416
$ws.Range("L7").Value = 4657.375
$ws.Range("M7").Value = -45457304
$ws.Range("N7").Value = -4881.375
$ws.Range("H16").Value = 787.8
$ws.Range("I16").Value = 835.2
$ws.Range("J16").Value = 740.4
$ws.Range("K16").Value = 835.2
$ws.Range("L16").Value = 740.4
$ws.Range("M16").Value = -665.2
$ws.Range("N16").Value = -1080.4
$ws.Range("H61").Value = 4694.3125
$ws.Range("I61").Value = 2738.6
$ws.Range("J61").Value = 11679
$ws.Range("K61").Value = 2738.6
$ws.Range("L61").Value = 11679
$ws.Range("M61").Value = -2536.6
$ws.Range("N61").Value = -12083
$ws.Range("H100").Value = 2279.087
$ws.Range("I100").Value = 2008.7142
$ws.Range("J100").Value = 2699.6667
$ws.Range("K100").Value = 2008.7142
$ws.Range("L100").Value = 2699.6667
$ws.Range("M100").Value = -1467.7142
$ws.Range("N100").Value = -3781.6667
$ws.Range("H113").Value = 4694.3125
$ws.Range("I113").Value = 2738.6
$ws.Range("J113").Value = 11679
$ws.Range("K113").Value = 2738.6
$ws.Range("L113").Value = 11679
$ws.Range("M113").Value = -568.5999999999999
$ws.Range("N113").Value = -16019
$ws.Range("H126").Value = 26319412
$ws.Range("I126").Value = 45457416
$ws.Range("J126").Value = 4657.375
$ws.Range("K126").Value = 136372248
$ws.Range("L126").Value = 13972.125
$ws.Range("M126").Value = -136369778
$ws.Range("N126").Value = -18912.125
$ws.Range("H132").Value = 7552.647
$ws.Range("I132").Value = 3800.8
$ws.Range("J132").Value = 9115.916999999999
$ws.Range("K132").Value = 11402.4
$ws.Range("L132").Value = 27347.751
$ws.Range("M132").Value = -8872.400000000001
$ws.Range("N132").Value = -32407.751

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 62506910
$ws.Range("I62").Value = 6999.6665
$ws.Range("J62").Value = 76929970
$ws.Range("K62").Value = 6999.6665
$ws.Range("L62").Value = 76929970
$ws.Range("M62").Value = -6375.6665
$ws.Range("N62").Value = -76931218
$ws.Range("H65").Value = 62506910
$ws.Range("I65").Value = 6999.6665
$ws.Range("J65").Value = 76929970
$ws.Range("K65").Value = 34998.3325
$ws.Range("L65").Value = 384649850
$ws.Range("M65").Value = -31878.3325
$ws.Range("N65").Value = -384656090
$ws.Range("H132").Value = 5500.1577
$ws.Range("I132").Value = 5321.7144
$ws.Range("J132").Value = 5999.8
$ws.Range("K132").Value = 15965.1432
$ws.Range("L132").Value = 17999.4
$ws.Range("M132").Value = -13435.1432
$ws.Range("N132").Value = -23059.4
$ws.Range("H136").Value = 5750.9
$ws.Range("I136").Value = 2127.5
$ws.Range("K136").Value = 6382.5
$ws.Range("M136").Value = -3832.5
